$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.557.00'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '1.914.81'
$ws.Range("E3").Value = '  +5.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5220'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3972'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09713'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.11%  '
$ws.Range("E10").Value = '  +4.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.06'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.545'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.24%  '
$ws.Range("D14").Value = '1.914.07'
$ws.Range("E14").Value = '  +5.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.591'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.92%  '
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001140'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06654'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.59%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.345'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.63%  '
$ws.Range("D23").Value = '28.662.84'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.703'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.41%  '
$ws.Range("D27").Value = '2.136.38'
$ws.Range("E27").Value = '  +5.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.762'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.639'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.887'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06772'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02436'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.262'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2226'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.41%  '
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6447'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.97%  '
$ws.Range("E43").Value = '  +1.31%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6096'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.779'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.283'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.041'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.71%  '
$ws.Range("E51").Value = '  +2.47%  '
